$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 6..31 (columns A..T), taken from the post-edit
# state of the workbook. Columns A,B,C,E..K are week-to-week constants that
# happen to be identical before/after for existing rows, and are simply
# (re)written here too so the 3 brand-new rows (29-31) get them as well.
$data = @(
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44970,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",50,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44970,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Segunda",30,2500,2500,2500,"`$/bandeja 2 kilos","Provincia de Diguillín",1250,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44971,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",30,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44187,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",80,2800,3000,2900,"`$/bandeja 2 kilos","Provincia de Linares",1450,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44187,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",65,1400,1500,1446,"`$/envase 1 kilo","Provincia de Diguillín",1446,1)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44942,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",60,2500,2500,2500,"`$/bandeja 2 kilos","Provincia de Diguillín",1250,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44949,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",60,2800,3000,2900,"`$/bandeja 2 kilos","Provincia de Diguillín",1450,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44174,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",150,3700,3800,3747,"`$/bandeja 2 kilos","Provincia de Linares",1874,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44953,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",30,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44965,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",50,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44596,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",120,2500,2700,2600,"`$/bandeja 2 kilos","Provincia de Linares",1300,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44952,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",30,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44967,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",50,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44967,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Segunda",30,2500,2500,2500,"`$/bandeja 2 kilos","Provincia de Diguillín",1250,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44594,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",120,2500,2800,2650,"`$/bandeja 2 kilos","Provincia de Linares",1325,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44181,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",65,3600,3800,3692,"`$/bandeja 2 kilos","Provincia de Diguillín",1846,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44181,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",80,1800,2000,1875,"`$/envase 1 kilo","Provincia de Diguillín",1875,1)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44966,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Segunda",30,2500,2500,2500,"`$/bandeja 2 kilos","Provincia de Diguillín",1250,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44932,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",60,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44960,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Segunda",60,2500,2500,2500,"`$/bandeja 2 kilos","Provincia de Diguillín",1250,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44931,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",100,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44951,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",60,2800,3000,2900,"`$/bandeja 2 kilos","Provincia de Diguillín",1450,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44935,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",50,3000,3000,3000,"`$/bandeja 2 kilos","Provincia de Diguillín",1500,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44540,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",240,3500,3800,3650,"`$/bandeja 2 kilos","Región del Maule",1825,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44944,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",60,2500,2500,2500,"`$/bandeja 2 kilos","Provincia de Diguillín",1250,2)
    ,@(7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44539,16,"Fruta",100101,"Berries",100101001,"Arándano (blue)","Sin especificar","Primera",200,3800,4000,3900,"`$/bandeja 2 kilos","Región del Maule",1950,2)
)

$startRow = 6
$dateCol = 4   # column D ("Fecha") carries the custom date style
$dateFormat = $ws.Cells.Item(2, $dateCol).NumberFormat()

$r = $startRow
foreach ($row in $data) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $ws.Cells.Item($r, $dateCol).NumberFormat = $dateFormat
    $r = $r + 1
}

Write-Host "Done. Last row written:" ($r - 1)
Write-Host "Dimension check D6:" $ws.Cells.Item(6, 4).Value()
Write-Host "Dimension check D31:" $ws.Cells.Item(31, 4).Value()
